$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 13; existing rows 13-32 shift down to 14-33
$ws.Rows.Item(13).Insert()

# Populate the new row 13 with the new record (values match surrounding rows' style)
$ws.Cells.Item(13, 1).Value = 1
$ws.Cells.Item(13, 2).Value = "Agrícola del Norte S.A. de Arica"
$ws.Cells.Item(13, 3).Value = "Arica y Parinacota"
$ws.Cells.Item(13, 4).Value = 45238
$ws.Cells.Item(13, 4).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Cells.Item(13, 5).Value = 15
$ws.Cells.Item(13, 6).Value = 100114007
$ws.Cells.Item(13, 7).Value = "Jengibre"
$ws.Cells.Item(13, 8).Value = "Sin especificar"
$ws.Cells.Item(13, 9).Value = "Primera"
$ws.Cells.Item(13, 10).Value = 250
$ws.Cells.Item(13, 11).Value = 16000
$ws.Cells.Item(13, 12).Value = 17000
$ws.Cells.Item(13, 13).Value = 16500
$ws.Cells.Item(13, 14).Value = "`$/caja 13 kilos"
$ws.Cells.Item(13, 15).Value = "Perú"
$ws.Cells.Item(13, 16).Value = 1269
$ws.Cells.Item(13, 17).Value = 13
$ws.Cells.Item(13, 18).Value = "Hortaliza"
